$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 corresponds to the c160a3b9... handoff/handback pair.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 16:41:26"
$wsZhCn.Range("G3").Value = "2016-01-11 16:42:41"

# de-de sheet: row 3 corresponds to the c160a3b9... handoff/handback pair.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 16:41:53"
$wsDeDe.Range("G3").Value = "2016-01-11 16:43:16"
